$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 223: new last row -- gets the special bottom-border style family (copied from the ORIGINAL last row, 199, before it is restyled below)
$ws.Range("A199:L199").Copy()
$ws.Range("A223:L223").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N199").Copy()
$ws.Range("N223").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A223").Value2 = 45569.24636287037
$ws.Range("B223").Value2 = "bagj11532@gmail.com"
$ws.Range("C223").Value2 = "체욱학과"
$ws.Range("D223").Value2 = 20244120
$ws.Range("E223").Value2 = "박준형"
$ws.Range("F223").Value2 = "엘리트 문자"
$ws.Range("G223").Value2 = "한글"
$ws.Range("H223").Value2 = "2개"
$ws.Range("I223").Value2 = 0.8
$ws.Range("J223").Value2 = "영국"
$ws.Range("K223").Value2 = "2배 정도 실직할 가능성이 높다"
$ws.Range("L223").Value2 = "Black"
$ws.Range("N223").Value2 = "헐, 반 밖에 안 남았네."

# Row 199 is no longer the last row; restyle it from special (21-24) to normal (16-19) family
$ws.Range("A195:L195").Copy()
$ws.Range("A199:L199").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N195").Copy()
$ws.Range("N199").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 200 (even/N-branch), template row 196
$ws.Range("A196:L196").Copy()
$ws.Range("A200:L200").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N196").Copy()
$ws.Range("N200").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A200").Value2 = 45568.76201637731
$ws.Range("B200").Value2 = "jenniferdy@naver.com"
$ws.Range("C200").Value2 = "간호학과"
$ws.Range("D200").Value2 = 20246206
$ws.Range("E200").Value2 = "권도연"
$ws.Range("F200").Value2 = "민주 문자"
$ws.Range("G200").Value2 = "한글"
$ws.Range("H200").Value2 = "1개"
$ws.Range("I200").Value2 = 0.8
$ws.Range("J200").Value2 = "미국"
$ws.Range("K200").Value2 = "사회활동이나 자원활동에 덜 참여한다"
$ws.Range("L200").Value2 = "Black"
$ws.Range("N200").Value2 = "헐, 반 밖에 안 남았네."

# Row 201 (odd/N-branch), template row 195
$ws.Range("A195:L195").Copy()
$ws.Range("A201:L201").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N195").Copy()
$ws.Range("N201").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A201").Value2 = 45568.76396357639
$ws.Range("B201").Value2 = "0214lily@naver.com"
$ws.Range("C201").Value2 = "체육학과"
$ws.Range("D201").Value2 = 20214104
$ws.Range("E201").Value2 = "김가희"
$ws.Range("F201").Value2 = "민주 문자"
$ws.Range("G201").Value2 = "한글"
$ws.Range("H201").Value2 = "하나도 없다"
$ws.Range("I201").Value2 = 0.8
$ws.Range("J201").Value2 = "대한민국"
$ws.Range("K201").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L201").Value2 = "Black"
$ws.Range("N201").Value2 = "헐, 반 밖에 안 남았네."

# Row 202 (even/N-branch), template row 196
$ws.Range("A196:L196").Copy()
$ws.Range("A202:L202").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N196").Copy()
$ws.Range("N202").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A202").Value2 = 45568.7716172338
$ws.Range("B202").Value2 = "eunse051013@naver.com"
$ws.Range("C202").Value2 = "간호학과"
$ws.Range("D202").Value2 = 20246225
$ws.Range("E202").Value2 = "김은세"
$ws.Range("F202").Value2 = "민주 문자"
$ws.Range("G202").Value2 = "한글"
$ws.Range("H202").Value2 = "4개"
$ws.Range("I202").Value2 = 0.8
$ws.Range("J202").Value2 = "대한민국"
$ws.Range("K202").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L202").Value2 = "Black"
$ws.Range("N202").Value2 = "휴우, 그래도 반이나 남았네."

# Row 203 (odd/M-branch), template row 197
$ws.Range("A197:L197").Copy()
$ws.Range("A203:L203").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M197").Copy()
$ws.Range("M203").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A203").Value2 = 45568.77559520833
$ws.Range("B203").Value2 = "hsjenny99@gmail.com"
$ws.Range("C203").Value2 = "소프트웨어학부"
$ws.Range("D203").Value2 = 20245246
$ws.Range("E203").Value2 = "전소현"
$ws.Range("F203").Value2 = "민주 문자"
$ws.Range("G203").Value2 = "한글"
$ws.Range("H203").Value2 = "1개"
$ws.Range("I203").Value2 = 0.8
$ws.Range("J203").Value2 = "대한민국"
$ws.Range("K203").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L203").Value2 = "Red"
$ws.Range("M203").Value2 = "휴우, 그래도 반이나 남았네."

# Row 204 (even/M-branch), template row 194
$ws.Range("A194:L194").Copy()
$ws.Range("A204:L204").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M194").Copy()
$ws.Range("M204").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A204").Value2 = 45568.7817197338
$ws.Range("B204").Value2 = "qudcksl1216@gmail.com"
$ws.Range("C204").Value2 = "경영"
$ws.Range("D204").Value2 = 20192926
$ws.Range("E204").Value2 = "윤병찬"
$ws.Range("F204").Value2 = "민주 문자"
$ws.Range("G204").Value2 = "한글"
$ws.Range("H204").Value2 = "하나도 없다"
$ws.Range("I204").Value2 = 0.9
$ws.Range("J204").Value2 = "영국"
$ws.Range("K204").Value2 = "건강이 좋지 않다"
$ws.Range("L204").Value2 = "Red"
$ws.Range("M204").Value2 = "휴우, 그래도 반이나 남았네."

# Row 205 (odd/M-branch), template row 197
$ws.Range("A197:L197").Copy()
$ws.Range("A205:L205").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M197").Copy()
$ws.Range("M205").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A205").Value2 = 45568.789298020834
$ws.Range("B205").Value2 = "0223wltn@naver.com"
$ws.Range("C205").Value2 = "식품영양학과"
$ws.Range("D205").Value2 = 20243850
$ws.Range("E205").Value2 = "홍지수"
$ws.Range("F205").Value2 = "민주 문자"
$ws.Range("G205").Value2 = "한글"
$ws.Range("H205").Value2 = "2개"
$ws.Range("I205").Value2 = 0.8
$ws.Range("J205").Value2 = "대한민국"
$ws.Range("K205").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L205").Value2 = "Red"
$ws.Range("M205").Value2 = "휴우, 그래도 반이나 남았네."

# Row 206 (even/M-branch), template row 194
$ws.Range("A194:L194").Copy()
$ws.Range("A206:L206").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M194").Copy()
$ws.Range("M206").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A206").Value2 = 45568.803086562504
$ws.Range("B206").Value2 = "bvc023@naver.com"
$ws.Range("C206").Value2 = "사회복지학과"
$ws.Range("D206").Value2 = 20217035
$ws.Range("E206").Value2 = "김수영"
$ws.Range("F206").Value2 = "민주 문자"
$ws.Range("G206").Value2 = "한글"
$ws.Range("H206").Value2 = "하나도 없다"
$ws.Range("I206").Value2 = 0.5
$ws.Range("J206").Value2 = "미국"
$ws.Range("K206").Value2 = "남들을 덜 신뢰한다"
$ws.Range("L206").Value2 = "Red"
$ws.Range("M206").Value2 = "헐, 반 밖에 안 남았네."

# Row 207 (odd/M-branch), template row 197
$ws.Range("A197:L197").Copy()
$ws.Range("A207:L207").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M197").Copy()
$ws.Range("M207").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A207").Value2 = 45568.8037778125
$ws.Range("B207").Value2 = "ckswo00@gmail.com"
$ws.Range("C207").Value2 = "빅데이터"
$ws.Range("D207").Value2 = 20217151
$ws.Range("E207").Value2 = "이찬재"
$ws.Range("F207").Value2 = "민주 문자"
$ws.Range("G207").Value2 = "한글"
$ws.Range("H207").Value2 = "2개"
$ws.Range("I207").Value2 = 0.2
$ws.Range("J207").Value2 = "대한민국"
$ws.Range("K207").Value2 = "건강이 좋지 않다"
$ws.Range("L207").Value2 = "Red"
$ws.Range("M207").Value2 = "휴우, 그래도 반이나 남았네."

# Row 208 (even/M-branch), template row 194
$ws.Range("A194:L194").Copy()
$ws.Range("A208:L208").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M194").Copy()
$ws.Range("M208").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A208").Value2 = 45568.86027855324
$ws.Range("B208").Value2 = "anfytlrtk3@naver.com"
$ws.Range("C208").Value2 = "식품영양학과"
$ws.Range("D208").Value2 = 20243830
$ws.Range("E208").Value2 = "윤가영"
$ws.Range("F208").Value2 = "민주 문자"
$ws.Range("G208").Value2 = "한글"
$ws.Range("H208").Value2 = "1개"
$ws.Range("I208").Value2 = 0.8
$ws.Range("J208").Value2 = "대한민국"
$ws.Range("K208").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L208").Value2 = "Red"
$ws.Range("M208").Value2 = "휴우, 그래도 반이나 남았네."

# Row 209 (odd/N-branch), template row 195
$ws.Range("A195:L195").Copy()
$ws.Range("A209:L209").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N195").Copy()
$ws.Range("N209").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A209").Value2 = 45568.866741481484
$ws.Range("B209").Value2 = "aldidrhemdgkrry1234@naver.com"
$ws.Range("C209").Value2 = "영어영문"
$ws.Range("D209").Value2 = 20241204
$ws.Range("E209").Value2 = "김용우"
$ws.Range("F209").Value2 = "민주 문자"
$ws.Range("G209").Value2 = "한글"
$ws.Range("H209").Value2 = "하나도 없다"
$ws.Range("I209").Value2 = 0.2
$ws.Range("J209").Value2 = "대한민국"
$ws.Range("K209").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L209").Value2 = "Black"
$ws.Range("N209").Value2 = "휴우, 그래도 반이나 남았네."

# Row 210 (even/N-branch), template row 196
$ws.Range("A196:L196").Copy()
$ws.Range("A210:L210").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N196").Copy()
$ws.Range("N210").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A210").Value2 = 45568.884467881944
$ws.Range("B210").Value2 = "sally200408@gmail.com"
$ws.Range("C210").Value2 = "경제학과"
$ws.Range("D210").Value2 = 20212818
$ws.Range("E210").Value2 = "박경화"
$ws.Range("F210").Value2 = "민주 문자"
$ws.Range("G210").Value2 = "한글"
$ws.Range("H210").Value2 = "하나도 없다"
$ws.Range("I210").Value2 = 0.8
$ws.Range("J210").Value2 = "대한민국"
$ws.Range("K210").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L210").Value2 = "Black"
$ws.Range("N210").Value2 = "헐, 반 밖에 안 남았네."

# Row 211 (odd/M-branch), template row 197
$ws.Range("A197:L197").Copy()
$ws.Range("A211:L211").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M197").Copy()
$ws.Range("M211").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A211").Value2 = 45568.893011631946
$ws.Range("B211").Value2 = "kyg031988@gmail.com"
$ws.Range("C211").Value2 = "일본학과"
$ws.Range("D211").Value2 = 20231606
$ws.Range("E211").Value2 = "김윤겸"
$ws.Range("F211").Value2 = "민주 문자"
$ws.Range("G211").Value2 = "한글"
$ws.Range("H211").Value2 = "하나도 없다"
$ws.Range("I211").Value2 = 0.8
$ws.Range("J211").Value2 = "대한민국"
$ws.Range("K211").Value2 = "사회활동이나 자원활동에 덜 참여한다"
$ws.Range("L211").Value2 = "Red"
$ws.Range("M211").Value2 = "휴우, 그래도 반이나 남았네."

# Row 212 (even/N-branch), template row 196
$ws.Range("A196:L196").Copy()
$ws.Range("A212:L212").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N196").Copy()
$ws.Range("N212").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A212").Value2 = 45568.93230721065
$ws.Range("B212").Value2 = "bagminhyeog534@gmail.com"
$ws.Range("C212").Value2 = "사회복지학과"
$ws.Range("D212").Value2 = 20242320
$ws.Range("E212").Value2 = "박민혁"
$ws.Range("F212").Value2 = "엘리트 문자"
$ws.Range("G212").Value2 = "한글"
$ws.Range("H212").Value2 = "하나도 없다"
$ws.Range("I212").Value2 = 0.1
$ws.Range("J212").Value2 = "이탈리아"
$ws.Range("K212").Value2 = "2배 정도 실직할 가능성이 높다"
$ws.Range("L212").Value2 = "Black"
$ws.Range("N212").Value2 = "헐, 반 밖에 안 남았네."

# Row 213 (odd/N-branch), template row 195
$ws.Range("A195:L195").Copy()
$ws.Range("A213:L213").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N195").Copy()
$ws.Range("N213").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A213").Value2 = 45568.99201777778
$ws.Range("B213").Value2 = "snp040609@naver.com"
$ws.Range("C213").Value2 = "경영학과"
$ws.Range("D213").Value2 = 20242957
$ws.Range("E213").Value2 = "박세나"
$ws.Range("F213").Value2 = "민주 문자"
$ws.Range("G213").Value2 = "한자"
$ws.Range("H213").Value2 = "2개"
$ws.Range("I213").Value2 = 0.8
$ws.Range("J213").Value2 = "대한민국"
$ws.Range("K213").Value2 = "2배 정도 실직할 가능성이 높다"
$ws.Range("L213").Value2 = "Black"
$ws.Range("N213").Value2 = "휴우, 그래도 반이나 남았네."

# Row 214 (even/M-branch), template row 194
$ws.Range("A194:L194").Copy()
$ws.Range("A214:L214").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M194").Copy()
$ws.Range("M214").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A214").Value2 = 45569.008761018515
$ws.Range("B214").Value2 = "wizkids0418@naver.con"
$ws.Range("C214").Value2 = "사회복지학부"
$ws.Range("D214").Value2 = 20242327
$ws.Range("E214").Value2 = "박혜인"
$ws.Range("F214").Value2 = "민주 문자"
$ws.Range("G214").Value2 = "한글"
$ws.Range("H214").Value2 = "1개"
$ws.Range("I214").Value2 = 0.8
$ws.Range("J214").Value2 = "대한민국"
$ws.Range("K214").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L214").Value2 = "Red"
$ws.Range("M214").Value2 = "헐, 반 밖에 안 남았네."

# Row 215 (odd/M-branch), template row 197
$ws.Range("A197:L197").Copy()
$ws.Range("A215:L215").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M197").Copy()
$ws.Range("M215").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A215").Value2 = 45569.025794166664
$ws.Range("B215").Value2 = "ilylive999@gmail.com"
$ws.Range("C215").Value2 = "일본학과"
$ws.Range("D215").Value2 = 20231621
$ws.Range("E215").Value2 = "이가현"
$ws.Range("F215").Value2 = "민주 문자"
$ws.Range("G215").Value2 = "한글"
$ws.Range("H215").Value2 = "하나도 없다"
$ws.Range("I215").Value2 = 0.8
$ws.Range("J215").Value2 = "대한민국"
$ws.Range("K215").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L215").Value2 = "Red"
$ws.Range("M215").Value2 = "휴우, 그래도 반이나 남았네."

# Row 216 (even/N-branch), template row 196
$ws.Range("A196:L196").Copy()
$ws.Range("A216:L216").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N196").Copy()
$ws.Range("N216").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A216").Value2 = 45569.033728229166
$ws.Range("B216").Value2 = "qwe92517@gmail.com"
$ws.Range("C216").Value2 = "간호학과"
$ws.Range("D216").Value2 = 20246271
$ws.Range("E216").Value2 = "이진영"
$ws.Range("F216").Value2 = "민주 문자"
$ws.Range("G216").Value2 = "한글"
$ws.Range("H216").Value2 = "하나도 없다"
$ws.Range("I216").Value2 = 0.5
$ws.Range("J216").Value2 = "대한민국"
$ws.Range("K216").Value2 = "건강이 좋지 않다"
$ws.Range("L216").Value2 = "Black"
$ws.Range("N216").Value2 = "헐, 반 밖에 안 남았네."

# Row 217 (odd/N-branch), template row 195
$ws.Range("A195:L195").Copy()
$ws.Range("A217:L217").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N195").Copy()
$ws.Range("N217").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A217").Value2 = 45569.05440775463
$ws.Range("B217").Value2 = "efgh124@naver.com"
$ws.Range("C217").Value2 = "간호"
$ws.Range("D217").Value2 = 20246282
$ws.Range("E217").Value2 = "정윤서"
$ws.Range("F217").Value2 = "민주 문자"
$ws.Range("G217").Value2 = "한글"
$ws.Range("H217").Value2 = "하나도 없다"
$ws.Range("I217").Value2 = 0.9
$ws.Range("J217").Value2 = "대한민국"
$ws.Range("K217").Value2 = "건강이 좋지 않다"
$ws.Range("L217").Value2 = "Black"
$ws.Range("N217").Value2 = "헐, 반 밖에 안 남았네."

# Row 218 (even/M-branch), template row 194
$ws.Range("A194:L194").Copy()
$ws.Range("A218:L218").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M194").Copy()
$ws.Range("M218").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A218").Value2 = 45569.06028276621
$ws.Range("B218").Value2 = "nyj7013@naver.com"
$ws.Range("C218").Value2 = "간호학과"
$ws.Range("D218").Value2 = 20246235
$ws.Range("E218").Value2 = "노예진"
$ws.Range("F218").Value2 = "민주 문자"
$ws.Range("G218").Value2 = "한글"
$ws.Range("H218").Value2 = "1개"
$ws.Range("I218").Value2 = 0.8
$ws.Range("J218").Value2 = "대한민국"
$ws.Range("K218").Value2 = "시간당 중위 임금이 60% 낮다"
$ws.Range("L218").Value2 = "Red"
$ws.Range("M218").Value2 = "휴우, 그래도 반이나 남았네."

# Row 219 (odd/N-branch), template row 195
$ws.Range("A195:L195").Copy()
$ws.Range("A219:L219").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N195").Copy()
$ws.Range("N219").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A219").Value2 = 45569.0604349537
$ws.Range("B219").Value2 = "yohihong@gmail.com"
$ws.Range("C219").Value2 = "소프트웨어학부"
$ws.Range("D219").Value2 = 20245161
$ws.Range("E219").Value2 = "민홍기"
$ws.Range("F219").Value2 = "민주 문자"
$ws.Range("G219").Value2 = "한글"
$ws.Range("H219").Value2 = "하나도 없다"
$ws.Range("I219").Value2 = 0.8
$ws.Range("J219").Value2 = "대한민국"
$ws.Range("K219").Value2 = "건강이 좋지 않다"
$ws.Range("L219").Value2 = "Black"
$ws.Range("N219").Value2 = "헐, 반 밖에 안 남았네."

# Row 220 (even/M-branch), template row 194
$ws.Range("A194:L194").Copy()
$ws.Range("A220:L220").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M194").Copy()
$ws.Range("M220").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A220").Value2 = 45569.079010914356
$ws.Range("B220").Value2 = "leejhzzang2005@naver.com"
$ws.Range("C220").Value2 = "사회복지학부"
$ws.Range("D220").Value2 = 20242342
$ws.Range("E220").Value2 = "이주현"
$ws.Range("F220").Value2 = "민주 문자"
$ws.Range("G220").Value2 = "한글"
$ws.Range("H220").Value2 = "하나도 없다"
$ws.Range("I220").Value2 = 0.8
$ws.Range("J220").Value2 = "영국"
$ws.Range("K220").Value2 = "건강이 좋지 않다"
$ws.Range("L220").Value2 = "Red"
$ws.Range("M220").Value2 = "휴우, 그래도 반이나 남았네."

# Row 221 (odd/M-branch), template row 197
$ws.Range("A197:L197").Copy()
$ws.Range("A221:L221").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("M197").Copy()
$ws.Range("M221").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A221").Value2 = 45569.0849347801
$ws.Range("B221").Value2 = "msy123581@naver.com"
$ws.Range("C221").Value2 = "중국학과"
$ws.Range("D221").Value2 = 20241519
$ws.Range("E221").Value2 = "문신영"
$ws.Range("F221").Value2 = "민주 문자"
$ws.Range("G221").Value2 = "한글"
$ws.Range("H221").Value2 = "3개"
$ws.Range("I221").Value2 = 0.8
$ws.Range("J221").Value2 = "대한민국"
$ws.Range("K221").Value2 = "남들을 덜 신뢰한다"
$ws.Range("L221").Value2 = "Red"
$ws.Range("M221").Value2 = "휴우, 그래도 반이나 남았네."

# Row 222 (even/N-branch), template row 196
$ws.Range("A196:L196").Copy()
$ws.Range("A222:L222").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N196").Copy()
$ws.Range("N222").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A222").Value2 = 45569.09024958333
$ws.Range("B222").Value2 = "minwl19@naver.com"
$ws.Range("C222").Value2 = "사회학과"
$ws.Range("D222").Value2 = 20217178
$ws.Range("E222").Value2 = "조민지"
$ws.Range("F222").Value2 = "민주 문자"
$ws.Range("G222").Value2 = "한글"
$ws.Range("H222").Value2 = "1개"
$ws.Range("I222").Value2 = 0.8
$ws.Range("J222").Value2 = "대한민국"
$ws.Range("K222").Value2 = "사회활동이나 자원활동에 덜 참여한다"
$ws.Range("L222").Value2 = "Black"
$ws.Range("N222").Value2 = "휴우, 그래도 반이나 남았네."

# Resize table to include new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N223"))
